# Update "想去人数" (want-to-go count) figures for several expo entries.
# Sheet "展览" holds the exhibition-only listing; sheet "全部类型" holds the
# combined listing of all event types. The same rows (by event) appear in
# both, so both must be updated to keep the data consistent.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 566
$wsExpo.Range("F5").Value = 406
$wsExpo.Range("F7").Value = 2380
$wsExpo.Range("F9").Value = 6102

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 566
$wsAll.Range("F5").Value = 406
$wsAll.Range("F9").Value = 2380
$wsAll.Range("F11").Value = 6102
